$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Sheet 2: "Monthly Trend" ---------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- New Sheet 3: "PO Forecast" -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Bring over the same formatting used on the other two sheets so the new
# sheet's header row / date column match the existing style indices.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A20").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$forecastData = @(
    @(45151.99999999999, 1, 0.6708354708721114, 1.91163704188749),
    @(45186.99999999999, 1, 0.7140810519330106, 1.973172314135023),
    @(45263.99999999999, 1, 0.8549604644081508, 2.109304780730848),
    @(45326.99999999999, 2, 0.9658123589031863, 2.191351050269728),
    @(45480.99999999999, 2, 1.229196513688165, 2.444875954282796),
    @(45564.99999999999, 2, 1.397847694180211, 2.624173951138218),
    @(45585.99999999999, 2, 1.427985924256478, 2.674916203390508),
    @(45599.99999999999, 2, 1.50548783124363, 2.659001994912831),
    @(45613.99999999999, 2, 1.507080641944797, 2.664895589591624),
    @(45634.99999999999, 2, 1.551946799721815, 2.730537143805542),
    @(45641.99999999999, 2, 1.510183815519915, 2.735702004937023),
    @(45648.99999999999, 2, 1.53917428474744, 2.756519614071765),
    @(45655.99999999999, 2, 1.528025032745859, 2.727470202511853),
    @(45662.99999999999, 2, 1.603426744644284, 2.780213667007029),
    @(45669.99999999999, 2, 1.57060283658408, 2.801023717721808),
    @(45676.99999999999, 2, 1.616430378592812, 2.783124582694605),
    @(45683.99999999999, 2, 1.61370539062034, 2.81864335307616),
    @(45690.99999999999, 2, 1.613988977162999, 2.842929909045568),
    @(45697.99999999999, 2, 1.658115533304002, 2.851511016780466)
)

for ($i = 0; $i -lt $forecastData.Length; $i++) {
    $r = $i + 2
    $row = $forecastData[$i]
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}
